# Apply cryptos list update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells that look numeric stay stored as text,
# matching the source sheet convention for the Price column.
$textForceCells = @(
    "D4", "D5", "D7", "D8", "D11", "D15", "D16", "D18",
    "D21", "D23", "D25", "D27", "D31", "D34", "D40", "D42",
    "D43", "D44", "D45", "D46", "D49", "D50"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (prices, volumes, and the two
# rank swaps: Polygon/Polkadot and PaxDollar/BitcoinSV).
$ws.Range('D2').Value = '29.186.10'
$ws.Range('E2').Value = '  +3.09%  '

$ws.Range('D3').Value = '1.580.51'
$ws.Range('E3').Value = '  +1.84%  '

$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -0.35%  '

$ws.Range('D5').Value = '212.39'
$ws.Range('E5').Value = '  +1.19%  '

$ws.Range('E6').Value = '  +6.71%  '

$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.36%  '

$ws.Range('D8').Value = '26.29'
$ws.Range('E8').Value = '  +10.10%  '

$ws.Range('E9').Value = '  +2.21%  '

$ws.Range('E10').Value = '  +1.80%  '

$ws.Range('D11').Value = '0.0904'
$ws.Range('E11').Value = '  +1.58%  '

$ws.Range('D12').Value = '1.806.06'

$ws.Range('D13').Value = '1.571.58'
$ws.Range('E13').Value = '  +1.31%  '

$ws.Range('D14').Value = '29.219.23'
$ws.Range('E14').Value = '  +3.16%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '3.71'
$ws.Range('E15').Value = '  +2.48%  '

$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '0.523'
$ws.Range('E16').Value = '  +2.53%  '

$ws.Range('E17').Value = '  +3.41%  '

$ws.Range('D18').Value = '237.24'
$ws.Range('E18').Value = '  +4.01%  '

$ws.Range('E19').Value = '  +1.73%  '

$ws.Range('E20').Value = '  +2.08%  '

$ws.Range('D21').Value = '0.997'
$ws.Range('E21').Value = '  -0.30%  '

$ws.Range('E22').Value = '  +1.76%  '

$ws.Range('D23').Value = '9.19'
$ws.Range('E23').Value = '  +2.88%  '

$ws.Range('E24').Value = '  +2.67%  '

$ws.Range('D25').Value = '154.41'
$ws.Range('E25').Value = '  +2.31%  '

$ws.Range('E26').Value = '  +5.15%  '

$ws.Range('D27').Value = '15.13'
$ws.Range('E27').Value = '  +2.56%  '

$ws.Range('E28').Value = '  +2.02%  '

$ws.Range('E29').Value = '  -0.29%  '

$ws.Range('E30').Value = '  +0.14%  '

$ws.Range('D31').Value = '1.06'
$ws.Range('E31').Value = '  +0.27%  '

$ws.Range('E32').Value = '  +1.70%  '

$ws.Range('D33').Value = '1.423.44'
$ws.Range('E33').Value = '  +2.48%  '

$ws.Range('D34').Value = '3.08'
$ws.Range('E34').Value = '  +1.87%  '

$ws.Range('E35').Value = '  -2.93%  '

$ws.Range('E36').Value = '  +1.94%  '

$ws.Range('E37').Value = '  +6.12%  '

$ws.Range('E38').Value = '  -1.78%  '

$ws.Range('E39').Value = '  +1.92%  '

$ws.Range('D40').Value = '0.532'
$ws.Range('E40').Value = '  +3.62%  '

$ws.Range('B42').Value = 'BitcoinSV'
$ws.Range('C42').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D42').Value = '52.98'
$ws.Range('E42').Value = '  +23.37%  '

$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = '0.997'
$ws.Range('E43').Value = '  -0.32%  '

$ws.Range('D44').Value = '0.789'
$ws.Range('E44').Value = '  +1.64%  '

$ws.Range('D45').Value = '0.0473'
$ws.Range('E45').Value = '  +3.58%  '

$ws.Range('D46').Value = '64.55'
$ws.Range('E46').Value = '  +4.29%  '

$ws.Range('E47').Value = '  -0.31%  '

$ws.Range('D48').Value = '1.717.61'
$ws.Range('E48').Value = '  +1.89%  '

$ws.Range('D49').Value = '0.836'
$ws.Range('E49').Value = '  -7.41%  '

$ws.Range('D50').Value = '85.28'
$ws.Range('E50').Value = '  -0.47%  '

$ws.Range('D51').Value = '0.0₆0101'
$ws.Range('E51').Value = '  -1.70%  '
